# Applies crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'62.025.79"
$ws.Range("E2").Value = "  +1.63%  "

# Row 3
$ws.Range("D3").Value = "'3.416.13"
$ws.Range("E3").Value = "  +0.87%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "'578.34"
$ws.Range("E5").Value = "  +1.31%  "

# Row 6
$ws.Range("D6").Value = "'144.76"
$ws.Range("E6").Value = "  +2.87%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.475"
$ws.Range("E8").Value = "  +0.27%  "

# Row 9
$ws.Range("E9").Value = "  -0.79%  "

# Row 10
$ws.Range("E10").Value = "  +0.90%  "

# Row 11
$ws.Range("D11").Value = "'0.386"
$ws.Range("E11").Value = "  -0.04%  "

# Row 12
$ws.Range("D12").Value = "'4.002.70"
$ws.Range("E12").Value = "  +0.96%  "

# Row 13
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'28.60"
$ws.Range("E13").Value = "  +2.66%  "

# Row 14
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "'0.125"
$ws.Range("E14").Value = "  -0.75%  "

# Row 15
$ws.Range("D15").Value = "'3.406.82"
$ws.Range("E15").Value = "  +0.25%  "

# Row 16
$ws.Range("E16").Value = "  +0.29%  "

# Row 17
$ws.Range("D17").Value = "'62.058.44"
$ws.Range("E17").Value = "  +1.51%  "

# Row 18
$ws.Range("D18").Value = "'6.17"
$ws.Range("E18").Value = "  +1.35%  "

# Row 19
$ws.Range("D19").Value = "'14.08"
$ws.Range("E19").Value = "  +3.50%  "

# Row 20
$ws.Range("D20").Value = "'9.18"

# Row 21
$ws.Range("D21").Value = "'390.85"
$ws.Range("E21").Value = "  +2.23%  "

# Row 22
$ws.Range("D22").Value = "'75.13"
$ws.Range("E22").Value = "  -0.71%  "

# Row 23
$ws.Range("D23").Value = "'0.555"
$ws.Range("E23").Value = "  +0.86%  "

# Row 25
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000115"
$ws.Range("E25").Value = "  +0.07%  "

# Row 26
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "'3.554.89"
$ws.Range("E26").Value = "  +0.94%  "

# Row 27
$ws.Range("D27").Value = "'0.186"
$ws.Range("E27").Value = "  -2.26%  "

# Row 28
$ws.Range("D28").Value = "'7.49"
$ws.Range("E28").Value = "  +3.35%  "

# Row 29
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("D30").Value = "'8.04"
$ws.Range("E30").Value = "  +0.70%  "

# Row 31
$ws.Range("E31").Value = "  +0.28%  "

# Row 32
$ws.Range("E32").Value = "  +0.11%  "

# Row 33
$ws.Range("E33").Value = "  +1.83%  "

# Row 34
$ws.Range("D34").Value = "'23.61"
$ws.Range("E34").Value = "  +1.39%  "

# Row 35
$ws.Range("D35").Value = "'5.31"
$ws.Range("E35").Value = "  +6.67%  "

# Row 36
$ws.Range("D36").Value = "'6.98"
$ws.Range("E36").Value = "  +0.64%  "

# Row 37
$ws.Range("D37").Value = "'168.20"
$ws.Range("E37").Value = "  +0.81%  "

# Row 38
$ws.Range("D38").Value = "'1.55"
$ws.Range("E38").Value = "  +6.30%  "

# Row 39
$ws.Range("D39").Value = "'3.449.57"
$ws.Range("E39").Value = "  +0.85%  "

# Row 40
$ws.Range("D40").Value = "'28.79"
$ws.Range("E40").Value = "  +9.02%  "

# Row 41
$ws.Range("D41").Value = "'0.0755"
$ws.Range("E41").Value = "  -1.43%  "

# Row 42
$ws.Range("D42").Value = "'0.786"
$ws.Range("E42").Value = "  +0.90%  "

# Row 43
$ws.Range("D43").Value = "'4.44"
$ws.Range("E43").Value = "  +1.79%  "

# Row 44
$ws.Range("D44").Value = "'1.68"
$ws.Range("E44").Value = "  +2.15%  "

# Row 45
$ws.Range("D45").Value = "'1.17"
$ws.Range("E45").Value = "  +4.49%  "

# Row 46
$ws.Range("D46").Value = "'2.504.33"
$ws.Range("E46").Value = "  +2.38%  "

# Row 47
$ws.Range("D47").Value = "'22.85"
$ws.Range("E47").Value = "  -0.43%  "

# Row 48
$ws.Range("D48").Value = "'6.65"
$ws.Range("E48").Value = "  +0.18%  "

# Row 49
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -0.09%  "

# Row 50
$ws.Range("D50").Value = "'0.0263"
$ws.Range("E50").Value = "  +0.76%  "

# Row 51
$ws.Range("D51").Value = "'2.09"
$ws.Range("E51").Value = "  -1.14%  "
